# Generate Report for Handoff
# - Updates the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it is shown (Overview!E2:F2, zh-cn!C2,
#   de-de!C2).
# - Bumps the "Latest HO Xliff Generate Date" / "Latest Handback DateTime"
#   timestamps for the Overview sheet and de-de sheet (they shared the same
#   value) and the zh-cn sheet's "Latest Handback DateTime".
# - Narrows the zh-cn/de-de Status column (and the matching Overview
#   zh-cn/de-de columns) now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$oldOverviewDate = "2016-09-02 09:08:18"
$newOverviewDate = "2016-09-02 09:09:05"

$oldZhDate = "2016-09-02 09:08:13"
$newZhDate = "2016-09-02 09:08:57"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text -----------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Timestamps --------------------------------------------------------
$wsOverview.Range("G2").Value = $newOverviewDate
$wsZhCn.Range("H2").Value = $newZhDate
$wsDeDe.Range("H2").Value = $newOverviewDate

# --- Column widths -------------------------------------------------------
# 29.9777047293527 (~character width 29.98) -> 17.2159881591797 (~character
# width 17.22) for the columns that hold the status text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.29
$wsOverview.Columns.Item(6).ColumnWidth = 16.29
$wsZhCn.Columns.Item(3).ColumnWidth = 16.29
$wsDeDe.Columns.Item(3).ColumnWidth = 16.29
